# Apply updated cryptocurrency market data (prices & 1h volume deltas).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.256.35"
$ws.Range("E2").Value = "  +7.67%  "
$ws.Range("D3").Value = "3.529.96"
$ws.Range("E3").Value = "  +11.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "191.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +12.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "553.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.74%  "
$ws.Range("D7").Value = "3.520.16"
$ws.Range("E7").Value = "  +10.69%  "
$ws.Range("E8").Value = "  +2.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("E10").Value = "  +5.52%  "
$ws.Range("E11").Value = "  +16.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.54%  "
$ws.Range("E13").Value = "  +7.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.97%  "
$ws.Range("D15").Value = "4.080.52"
$ws.Range("E15").Value = "  +11.07%  "
$ws.Range("D16").Value = "3.525.45"
$ws.Range("E16").Value = "  +11.52%  "
$ws.Range("E17").Value = "  +3.70%  "
$ws.Range("D18").Value = "67.286.45"
$ws.Range("E18").Value = "  +8.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "429.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +18.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("E27").Value = "  +11.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "650.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.04%  "
$ws.Range("E32").Value = "  +4.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.16%  "
$ws.Range("E34").Value = "  +5.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "59.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "38.81"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.25%  "
$ws.Range("D37").Value = "0.0₃0821"
$ws.Range("E37").Value = "  +16.97%  "
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.391"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.86%  "
$ws.Range("E40").Value = "  +15.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +14.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "3.031.91"
$ws.Range("E43").Value = "  +5.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.72%  "
$ws.Range("E45").Value = "  +12.19%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.93%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +12.98%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0419"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.07%  "
$ws.Range("E49").Value = "  +6.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +15.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "140.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.95%  "
